$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": update scraped_at (column K) timestamps for rows 2-44 ---
$ws1 = $wb.Worksheets.Item("snapshot")

$timestamps = @(
    "2025-12-14T03:02:51.001682+00:00",
    "2025-12-14T03:02:51.001720+00:00",
    "2025-12-14T03:02:51.001740+00:00",
    "2025-12-14T03:02:53.611182+00:00",
    "2025-12-14T03:02:53.611208+00:00",
    "2025-12-14T03:02:53.611225+00:00",
    "2025-12-14T03:02:56.370525+00:00",
    "2025-12-14T03:02:58.992734+00:00",
    "2025-12-14T03:03:01.350770+00:00",
    "2025-12-14T03:03:03.637568+00:00",
    "2025-12-14T03:03:09.120506+00:00",
    "2025-12-14T03:03:09.120539+00:00",
    "2025-12-14T03:03:11.868796+00:00",
    "2025-12-14T03:03:14.532440+00:00",
    "2025-12-14T03:03:17.344489+00:00",
    "2025-12-14T03:03:20.024868+00:00",
    "2025-12-14T03:03:20.024896+00:00",
    "2025-12-14T03:03:20.024913+00:00",
    "2025-12-14T03:03:20.024929+00:00",
    "2025-12-14T03:03:20.024946+00:00",
    "2025-12-14T03:03:22.683614+00:00",
    "2025-12-14T03:03:22.683645+00:00",
    "2025-12-14T03:03:22.683664+00:00",
    "2025-12-14T03:03:25.436484+00:00",
    "2025-12-14T03:03:25.436522+00:00",
    "2025-12-14T03:03:25.436546+00:00",
    "2025-12-14T03:03:25.436570+00:00",
    "2025-12-14T03:03:25.436592+00:00",
    "2025-12-14T03:03:28.216006+00:00",
    "2025-12-14T03:03:28.216034+00:00",
    "2025-12-14T03:03:28.216052+00:00",
    "2025-12-14T03:03:30.967610+00:00",
    "2025-12-14T03:03:30.967640+00:00",
    "2025-12-14T03:03:30.967659+00:00",
    "2025-12-14T03:03:33.702789+00:00",
    "2025-12-14T03:03:36.432702+00:00",
    "2025-12-14T03:03:36.432730+00:00",
    "2025-12-14T03:03:36.432746+00:00",
    "2025-12-14T03:03:39.281673+00:00",
    "2025-12-14T03:03:42.180973+00:00",
    "2025-12-14T03:03:42.181005+00:00",
    "2025-12-14T03:03:44.915136+00:00",
    "2025-12-14T03:03:44.915167+00:00"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 11).Value = $timestamps[$i]
}

# --- Sheet "returned": remove the single data row (row 2) ---
$ws2 = $wb.Worksheets.Item("returned")
$ws2.Rows.Item(2).Delete()

$wb.Save()
